# Updated cryptos list values (Price and Volume(1h) columns) to reflect
# refreshed market data, as produced by the scheduled GitHub Actions run.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A few refreshed Price values are plain decimal numbers (single decimal
# point). Force those specific cells to Text format first so Excel keeps
# storing/display them as exact strings (matching the sheet's existing
# inline-string convention) instead of silently re-typing them as floats.
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"

$ws.Range("D2").Value = '27.960.31'
$ws.Range("E2").Value = '  -0.36%  '
$ws.Range("D3").Value = '1.634.88'
$ws.Range("E3").Value = '  -1.09%  '
$ws.Range("E4").Value = '  +0.41%  '
$ws.Range("D5").Value = '212.07'
$ws.Range("E5").Value = '  -0.90%  '
$ws.Range("E6").Value = '  -0.96%  '
$ws.Range("E7").Value = '  +0.30%  '
$ws.Range("E8").Value = '  -1.29%  '
$ws.Range("E9").Value = '  -2.91%  '
$ws.Range("E10").Value = '  -0.08%  '
$ws.Range("D11").Value = '0.0883'
$ws.Range("E11").Value = '  +1.10%  '
$ws.Range("D12").Value = '1.866.29'
$ws.Range("E12").Value = '  -1.00%  '
$ws.Range("D13").Value = '1.635.02'
$ws.Range("E13").Value = '  -1.11%  '
$ws.Range("E14").Value = '  -0.52%  '
$ws.Range("E15").Value = '  -0.30%  '
$ws.Range("E16").Value = '  -0.55%  '
$ws.Range("D17").Value = '27.960.86'
$ws.Range("E17").Value = '  -0.28%  '
$ws.Range("D18").Value = '230.82'
$ws.Range("E18").Value = '  -1.09%  '
$ws.Range("E19").Value = '  -0.07%  '
$ws.Range("D20").Value = '7.54'
$ws.Range("E20").Value = '  -2.38%  '
$ws.Range("E21").Value = '  +0.19%  '
$ws.Range("E22").Value = '  -0.67%  '
$ws.Range("D23").Value = '10.37'
$ws.Range("E23").Value = '  -3.10%  '
$ws.Range("E24").Value = '  -4.01%  '
$ws.Range("D25").Value = '154.96'
$ws.Range("E25").Value = '  +1.83%  '
$ws.Range("E26").Value = '  +0.62%  '
$ws.Range("D27").Value = '15.67'
$ws.Range("E27").Value = '  -0.75%  '
$ws.Range("E28").Value = '  -0.68%  '
$ws.Range("E29").Value = '  +0.29%  '
$ws.Range("E30").Value = '  -0.90%  '
$ws.Range("E31").Value = '  -0.25%  '
$ws.Range("E32").Value = '  +0.96%  '
$ws.Range("D33").Value = '1.407.25'
$ws.Range("E33").Value = '  -3.11%  '
$ws.Range("E34").Value = '  -0.38%  '
$ws.Range("E35").Value = '  +0.05%  '
$ws.Range("E36").Value = '  +9.41%  '
$ws.Range("E37").Value = '  +1.77%  '
$ws.Range("E38").Value = '  +0.40%  '
$ws.Range("D39").Value = '0.562'
$ws.Range("E39").Value = '  +0.03%  '
$ws.Range("E40").Value = '  -2.18%  '
$ws.Range("E41").Value = '  +0.14%  '
$ws.Range("E42").Value = '  +0.12%  '
$ws.Range("E43").Value = '  -3.75%  '
$ws.Range("E45").Value = '  -0.24%  '
$ws.Range("E46").Value = '  -1.34%  '
$ws.Range("D47").Value = '1.775.88'
$ws.Range("E47").Value = '  -0.94%  '
$ws.Range("D48").Value = '87.92'
$ws.Range("E48").Value = '  -1.20%  '
$ws.Range("E49").Value = '  +6.45%  '
$ws.Range("E50").Value = '  -1.20%  '
$ws.Range("E51").Value = '  -0.27%  '
